$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-11 20:30:31"
$wsZhCn.Range("H3").Value = "2016-03-11 20:30:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-11 20:30:34"
$wsDeDe.Range("H3").Value = "2016-03-11 20:30:55"
